$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 35
$ws1.Range("F3").Value = 94
$ws1.Range("F4").Value = 1435
$ws1.Range("F5").Value = 176
$ws1.Range("F6").Value = 35
$ws1.Range("F8").Value = 9647
$ws1.Range("F10").Value = 107
$ws1.Range("F11").Value = 239
$ws1.Range("F13").Value = 365
$ws1.Range("F14").Value = 6646
$ws1.Range("F15").Value = 1079
$ws1.Range("F16").Value = 120
$ws1.Range("F17").Value = 48
$ws1.Range("F18").Value = 177

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 35
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 1435
$ws4.Range("F5").Value = 176
$ws4.Range("F6").Value = 35
$ws4.Range("F10").Value = 9647
$ws4.Range("F12").Value = 107
$ws4.Range("F13").Value = 239
$ws4.Range("F15").Value = 365
$ws4.Range("F16").Value = 6646
$ws4.Range("F17").Value = 1079
$ws4.Range("F18").Value = 120
$ws4.Range("F19").Value = 48
$ws4.Range("F20").Value = 177

$wb.Save()
